$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly entry: 星期三 (Wed) 2019-04-24, Hibernate study + star project build.
$ws.Range("A26").Value = "2019年4月24日23:01:24"
$ws.Range("B26").Value = "周三"
$ws.Range("C26").Value = "Hibernate练习"
$ws.Range("D26").Value = "8:30--10:10"

$ws.Range("C27").Value = "star项目构建（base，dao）"
$ws.Range("D27").Value = "15:00--15:40"

$ws.Range("C28").Value = "star项目构建（base，dao）"
$ws.Range("D28").Value = "23:00--24:00"

# Match the author's updated selection (cursor moved to D30 after the new rows).
$ws.Range("D30").Select()
